$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2 through 74 all hold the date serial 45186
# (2023-09-17). Update them to 45188 (2023-09-19).
for ($row = 2; $row -le 74; $row++) {
    $ws.Cells.Item($row, 3).Value = 45188
}
